$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range('J2').Value = 1435
$ws.Range('I3').Value = 7488
$ws.Range('J3').Value = 1501
$ws.Range('E4').Value = 1983
$ws.Range('I4').Value = 1754
$ws.Range('J4').Value = 343
$ws.Range('E6').Value = 9680
$ws.Range('J6').Value = 1961
$ws.Range('E7').Value = 25987
$ws.Range('J7').Value = 5345

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range('J2').Value = 22
$ws.Range('J7').Value = 62

$ws = $wb.Worksheets.Item('Fuller Park')
$ws.Range('J2').Value = 8
$ws.Range('J7').Value = 20

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range('J2').Value = 52
$ws.Range('J3').Value = 64
$ws.Range('J7').Value = 187

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range('J2').Value = 19
$ws.Range('J7').Value = 67

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range('J6').Value = 6
$ws.Range('J7').Value = 36

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range('J4').Value = 26
$ws.Range('J6').Value = 52
$ws.Range('J7').Value = 151
$ws.Range('J8').Value = 328
$ws.Range('J11').Value = 68
$ws.Range('J19').Value = 184
$ws.Range('J30').Value = 20
$ws.Range('J31').Value = 36
$ws.Range('E33').Value = 1508
$ws.Range('J37').Value = 187
$ws.Range('J42').Value = 212
$ws.Range('I43').Value = 224
$ws.Range('J47').Value = 46
$ws.Range('J50').Value = 27
$ws.Range('J51').Value = 70
$ws.Range('J52').Value = 115
$ws.Range('J53').Value = 50
$ws.Range('J54').Value = 107
$ws.Range('J57').Value = 22
$ws.Range('J61').Value = 10
$ws.Range('E63').Value = 329
$ws.Range('I63').Value = 191
$ws.Range('J63').Value = 23
$ws.Range('J66').Value = 11
$ws.Range('J73').Value = 52
$ws.Range('J76').Value = 87
$ws.Range('J79').Value = 169
$ws.Range('I83').Value = 566
$ws.Range('J83').Value = 132
$ws.Range('J85').Value = 243
$ws.Range('J88').Value = 52
$ws.Range('J89').Value = 62
$ws.Range('J91').Value = 69
$ws.Range('J94').Value = 40
$ws.Range('J97').Value = 35
$ws.Range('J99').Value = 67
$ws.Range('E101').Value = 25987
$ws.Range('J101').Value = 5345

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range('J2').Value = 44
$ws.Range('I5').Value = 24
$ws.Range('I7').Value = 566
$ws.Range('J7').Value = 132

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range('E6').Value = 608
$ws.Range('E7').Value = 1508

$ws = $wb.Worksheets.Item('Loop')
$ws.Range('J3').Value = 18
$ws.Range('J4').Value = 6
$ws.Range('J7').Value = 107

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range('J2').Value = 45
$ws.Range('J3').Value = 48
$ws.Range('J7').Value = 184

$ws = $wb.Worksheets.Item('River North')
$ws.Range('J3').Value = 20
$ws.Range('J6').Value = 51
$ws.Range('J7').Value = 87

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range('J2').Value = 59
$ws.Range('J3').Value = 95
$ws.Range('J7').Value = 243

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range('J2').Value = 17
$ws.Range('J6').Value = 18
$ws.Range('J7').Value = 52

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range('J6').Value = 116
$ws.Range('J7').Value = 212

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range('J4').Value = 6
$ws.Range('J7').Value = 69

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range('J2').Value = 44
$ws.Range('J6').Value = 49
$ws.Range('J7').Value = 169

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range('J2').Value = 28
$ws.Range('J7').Value = 115

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range('J6').Value = 23
$ws.Range('J7').Value = 40

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range('J6').Value = 24
$ws.Range('J7').Value = 46

$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Range('J2').Value = 9
$ws.Range('J4').Value = 5
$ws.Range('J7').Value = 27

$ws = $wb.Worksheets.Item('North Center')
$ws.Range('J6').Value = 6
$ws.Range('J7').Value = 11

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range('J2').Value = 22
$ws.Range('J7').Value = 68

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range('J4').Value = 7
$ws.Range('J6').Value = 9
$ws.Range('J7').Value = 52

$ws = $wb.Worksheets.Item('West Town')
$ws.Range('J6').Value = 23
$ws.Range('J7').Value = 35

$ws = $wb.Worksheets.Item('United Center')
$ws.Range('J2').Value = 9
$ws.Range('J6').Value = 27
$ws.Range('J7').Value = 52

$ws = $wb.Worksheets.Item('Austin')
$ws.Range('J2').Value = 108
$ws.Range('J3').Value = 108
$ws.Range('J4').Value = 17
$ws.Range('J7').Value = 328

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range('J2').Value = 17
$ws.Range('J7').Value = 70

$ws = $wb.Worksheets.Item('Mckinley Park')
$ws.Range('J4').Value = 2
$ws.Range('J7').Value = 22

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range('I3').Value = 38
$ws.Range('I7').Value = 224

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range('J6').Value = 28
$ws.Range('J7').Value = 50

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range('J3').Value = 48
$ws.Range('J6').Value = 49
$ws.Range('J7').Value = 151

$ws = $wb.Worksheets.Item('Archer Heights')
$ws.Range('J3').Value = 7
$ws.Range('J7').Value = 26

$ws = $wb.Worksheets.Item('Mount Greenwood')
$ws.Range('J4').Value = 3
$ws.Range('J7').Value = 10

Write-Host "Applied 127 cell updates"